$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = "CNReplication.setReplicationPolicy()"
$ws.Range("E9").Value = "manual (Tier 1), MNAuthorization.setAccessPolicy(), MNStorage.update ()(all must call CNCore.systemMetadataChanged())"

$ws.Range("E23").Select()
